# "nuevo cómputo de la tabla 2"
# Table 2 was recomputed; a handful of cells in columns H/I/J (rows 8-20)
# change to their newly-computed values. The row/column headers (column A
# and row 1) keep the same text-format styling they already had.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply the text format to the label column / header row so the whole
# table reads as freshly (re)written, matching the new computation pass.
$ws.Range("A1:M1").NumberFormat = "@"
$ws.Range("A1:A20").NumberFormat = "@"

# --- Row 8 ---
$ws.Range("H8").Value = 0.027

# --- Row 9 ---
$ws.Range("H9").Value = 0.97299999999999998

# --- Row 10 ---
$ws.Range("I10").Value = 0.024

# --- Row 11 ---
$ws.Range("I11").Value = 0.031
$ws.Range("J11").Value = 0.11600000000000001

# --- Row 12 ---
$ws.Range("I12").Value = 0.97599999999999998

# --- Row 13 ---
$ws.Range("I13").Value = 0.96899999999999997

# --- Row 17 ---
$ws.Range("H17").Value = -1009688.6090000001
$ws.Range("I17").Value = -1009688.612
$ws.Range("J17").Value = -499773.44300000003

# --- Row 18 ---
$ws.Range("H18").Value = -983529.43299999996
$ws.Range("I18").Value = -983529.43599999999
$ws.Range("J18").Value = -489291.13299999997

# --- Row 19 ---
$ws.Range("H19").Value = 1968630.6140000001
$ws.Range("I19").Value = 1968630.3060000001
$ws.Range("J19").Value = 963236.42099999997

# --- Row 20 ---
$ws.Range("H20").Value = -24587.429
$ws.Range("I20").Value = -24587.741000000002
$ws.Range("J20").Value = -25828.155999999999
